$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Generate Report for Handback
#
# Both localized files (e6081861-...md and 7b09e590-...md) have now been
# handed back ("Handed back: in sync with en-US") with fresh handback
# timestamps. The report rows were regenerated/re-sorted so the
# 7b09e590 file now appears first (row 2) and the e6081861 file second
# (row 3) on every sheet - but the hyperlink relationship ids stay bound
# to their original row position, so we only change each hyperlink's
# display text (and underlying cell value), never its target address.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "7b09e590-e9f9-4ec8-82b0-57ecea581437.md"
$ws1.Range("B2").Value = "Handed back: in sync with en-US"
$ws1.Range("C2").Value = "Handed back: in sync with en-US"
$ws1.Range("D2").Value = "2016-03-25 09:02:57"

$ws1.Range("A3").Value = "e6081861-f56e-4a77-9236-9901c9cdc7e9.md"
$ws1.Range("B3").Value = "Handed back: in sync with en-US"
$ws1.Range("C3").Value = "Handed back: in sync with en-US"
$ws1.Range("D3").Value = "2016-03-25 09:02:57"

# Rebuild hyperlinks, keeping each target Address pinned to its original
# row/ cell position (rId2 -> A2, rId3 -> A3) but refreshing the display
# text to the file that now lives in that row.
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/631b1e1a4a8696c944d547fa3611d0ae79356dba/e2e/e6081861-f56e-4a77-9236-9901c9cdc7e9.md", [Type]::Missing, [Type]::Missing, "7b09e590-e9f9-4ec8-82b0-57ecea581437.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/41cbf513a721ef7c2ad1032485e7843b5e96492e/e2e/7b09e590-e9f9-4ec8-82b0-57ecea581437.md", [Type]::Missing, [Type]::Missing, "e6081861-f56e-4a77-9236-9901c9cdc7e9.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "7b09e590-e9f9-4ec8-82b0-57ecea581437.md"
$ws2.Range("B2").Value = ".md"
$ws2.Range("C2").Value = "Handed back: in sync with en-US"
$ws2.Range("D2").Value = "7b09e590-e9f9-4ec8-82b0-57ecea581437.4fdc4ad2453416e45d2658ced24248438eeb5397.zh-cn.xlf"
$ws2.Range("E2").Value = "2016-03-25 09:02:48"
$ws2.Range("F2").Value = "7b09e590-e9f9-4ec8-82b0-57ecea581437.md"
$ws2.Range("G2").Value = "7b09e590-e9f9-4ec8-82b0-57ecea581437.4fdc4ad2453416e45d2658ced24248438eeb5397.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-03-25 09:03:38"
$ws2.Range("J2").Value = "Include"

$ws2.Range("A3").Value = "e6081861-f56e-4a77-9236-9901c9cdc7e9.md"
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Handed back: in sync with en-US"
$ws2.Range("D3").Value = "e6081861-f56e-4a77-9236-9901c9cdc7e9.65b374d3eee2adffd0001118f274db0000a84dd7.zh-cn.xlf"
$ws2.Range("E3").Value = "2016-03-25 09:02:48"
$ws2.Range("F3").Value = "e6081861-f56e-4a77-9236-9901c9cdc7e9.md"
$ws2.Range("G3").Value = "e6081861-f56e-4a77-9236-9901c9cdc7e9.65b374d3eee2adffd0001118f274db0000a84dd7.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-03-25 09:03:38"
$ws2.Range("J3").Value = "Include"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/631b1e1a4a8696c944d547fa3611d0ae79356dba/e2e/e6081861-f56e-4a77-9236-9901c9cdc7e9.md", [Type]::Missing, [Type]::Missing, "7b09e590-e9f9-4ec8-82b0-57ecea581437.md")
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9fc40aee9f66dafb6a0e8689701d56e999717828/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/mt/e6081861-f56e-4a77-9236-9901c9cdc7e9.65b374d3eee2adffd0001118f274db0000a84dd7.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "7b09e590-e9f9-4ec8-82b0-57ecea581437.4fdc4ad2453416e45d2658ced24248438eeb5397.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/aa3d1ef1149f9e28b38e91754c25acaa7cbe2fce/e2e/e6081861-f56e-4a77-9236-9901c9cdc7e9.md", [Type]::Missing, [Type]::Missing, "7b09e590-e9f9-4ec8-82b0-57ecea581437.md")
$ws2.Hyperlinks.Add($ws2.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c39556046412838b3b466374b8544c03231e43e2/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/e6081861-f56e-4a77-9236-9901c9cdc7e9.65b374d3eee2adffd0001118f274db0000a84dd7.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "7b09e590-e9f9-4ec8-82b0-57ecea581437.4fdc4ad2453416e45d2658ced24248438eeb5397.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/41cbf513a721ef7c2ad1032485e7843b5e96492e/e2e/7b09e590-e9f9-4ec8-82b0-57ecea581437.md", [Type]::Missing, [Type]::Missing, "e6081861-f56e-4a77-9236-9901c9cdc7e9.md")
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9fc40aee9f66dafb6a0e8689701d56e999717828/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/mt/7b09e590-e9f9-4ec8-82b0-57ecea581437.4fdc4ad2453416e45d2658ced24248438eeb5397.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "e6081861-f56e-4a77-9236-9901c9cdc7e9.65b374d3eee2adffd0001118f274db0000a84dd7.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/aa3d1ef1149f9e28b38e91754c25acaa7cbe2fce/e2e/7b09e590-e9f9-4ec8-82b0-57ecea581437.md", [Type]::Missing, [Type]::Missing, "e6081861-f56e-4a77-9236-9901c9cdc7e9.md")
$ws2.Hyperlinks.Add($ws2.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c39556046412838b3b466374b8544c03231e43e2/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/7b09e590-e9f9-4ec8-82b0-57ecea581437.4fdc4ad2453416e45d2658ced24248438eeb5397.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "e6081861-f56e-4a77-9236-9901c9cdc7e9.65b374d3eee2adffd0001118f274db0000a84dd7.zh-cn.xlf")

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "7b09e590-e9f9-4ec8-82b0-57ecea581437.md"
$ws3.Range("B2").Value = ".md"
$ws3.Range("C2").Value = "Handed back: in sync with en-US"
$ws3.Range("D2").Value = "7b09e590-e9f9-4ec8-82b0-57ecea581437.4fdc4ad2453416e45d2658ced24248438eeb5397.de-de.xlf"
$ws3.Range("E2").Value = "2016-03-25 09:02:57"
$ws3.Range("F2").Value = "7b09e590-e9f9-4ec8-82b0-57ecea581437.md"
$ws3.Range("G2").Value = "7b09e590-e9f9-4ec8-82b0-57ecea581437.4fdc4ad2453416e45d2658ced24248438eeb5397.de-de.xlf"
$ws3.Range("H2").Value = "2016-03-25 09:03:57"
$ws3.Range("J2").Value = "Include"

$ws3.Range("A3").Value = "e6081861-f56e-4a77-9236-9901c9cdc7e9.md"
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Handed back: in sync with en-US"
$ws3.Range("D3").Value = "e6081861-f56e-4a77-9236-9901c9cdc7e9.65b374d3eee2adffd0001118f274db0000a84dd7.de-de.xlf"
$ws3.Range("E3").Value = "2016-03-25 09:02:57"
$ws3.Range("F3").Value = "e6081861-f56e-4a77-9236-9901c9cdc7e9.md"
$ws3.Range("G3").Value = "e6081861-f56e-4a77-9236-9901c9cdc7e9.65b374d3eee2adffd0001118f274db0000a84dd7.de-de.xlf"
$ws3.Range("H3").Value = "2016-03-25 09:03:57"
$ws3.Range("J3").Value = "Include"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/631b1e1a4a8696c944d547fa3611d0ae79356dba/e2e/e6081861-f56e-4a77-9236-9901c9cdc7e9.md", [Type]::Missing, [Type]::Missing, "7b09e590-e9f9-4ec8-82b0-57ecea581437.md")
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/16a0860bf93499805bc619c858f1989b928f4300/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/mt/e6081861-f56e-4a77-9236-9901c9cdc7e9.65b374d3eee2adffd0001118f274db0000a84dd7.de-de.xlf", [Type]::Missing, [Type]::Missing, "7b09e590-e9f9-4ec8-82b0-57ecea581437.4fdc4ad2453416e45d2658ced24248438eeb5397.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/429d84775198908804fceb3ec2a1a6576ef428cb/e2e/e6081861-f56e-4a77-9236-9901c9cdc7e9.md", [Type]::Missing, [Type]::Missing, "7b09e590-e9f9-4ec8-82b0-57ecea581437.md")
$ws3.Hyperlinks.Add($ws3.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c5c53f544066256e0878d6906839cd4f6d5191ef/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/e6081861-f56e-4a77-9236-9901c9cdc7e9.65b374d3eee2adffd0001118f274db0000a84dd7.de-de.xlf", [Type]::Missing, [Type]::Missing, "7b09e590-e9f9-4ec8-82b0-57ecea581437.4fdc4ad2453416e45d2658ced24248438eeb5397.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/41cbf513a721ef7c2ad1032485e7843b5e96492e/e2e/7b09e590-e9f9-4ec8-82b0-57ecea581437.md", [Type]::Missing, [Type]::Missing, "e6081861-f56e-4a77-9236-9901c9cdc7e9.md")
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/16a0860bf93499805bc619c858f1989b928f4300/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/mt/7b09e590-e9f9-4ec8-82b0-57ecea581437.4fdc4ad2453416e45d2658ced24248438eeb5397.de-de.xlf", [Type]::Missing, [Type]::Missing, "e6081861-f56e-4a77-9236-9901c9cdc7e9.65b374d3eee2adffd0001118f274db0000a84dd7.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/429d84775198908804fceb3ec2a1a6576ef428cb/e2e/7b09e590-e9f9-4ec8-82b0-57ecea581437.md", [Type]::Missing, [Type]::Missing, "e6081861-f56e-4a77-9236-9901c9cdc7e9.md")
$ws3.Hyperlinks.Add($ws3.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c5c53f544066256e0878d6906839cd4f6d5191ef/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/7b09e590-e9f9-4ec8-82b0-57ecea581437.4fdc4ad2453416e45d2658ced24248438eeb5397.de-de.xlf", [Type]::Missing, [Type]::Missing, "e6081861-f56e-4a77-9236-9901c9cdc7e9.65b374d3eee2adffd0001118f274db0000a84dd7.de-de.xlf")
